# Refresh cached market-board price/profit figures across the per-job
# leve-profit sheets (scheduled data-refresh run).
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 71436320
$ws.Range("J64").Value = 8873.5
$ws.Range("L64").Value = 8873.5
$ws.Range("N64").Value = -9369.5
# Row 67
$ws.Range("H67").Value = 71436320
$ws.Range("J67").Value = 8873.5
$ws.Range("L67").Value = 8873.5
$ws.Range("N67").Value = -10589.5
# Row 76
$ws.Range("H76").Value = 19998.5
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 19998.5
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 19998.5
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -20628.5
# Row 79
$ws.Range("H79").Value = 19998.5
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 19998.5
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 19998.5
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -22182.5
# Row 87
$ws.Range("H87").Value = 54856.57
$ws.Range("J87").Value = 54856.57
$ws.Range("L87").Value = 54856.57
$ws.Range("N87").Value = -57352.57
# Row 90
$ws.Range("H90").Value = 54856.57
$ws.Range("J90").Value = 54856.57
$ws.Range("L90").Value = 164569.71
$ws.Range("N90").Value = -177049.71
# Row 98
$ws.Range("H98").Value = 1982.9756
$ws.Range("I98").Value = 1982.9756
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1982.9756
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -484.9756
$ws.Range("N98").ClearContents()
# Row 112
$ws.Range("H112").Value = 4937.8696
$ws.Range("J112").Value = 5514.7
$ws.Range("L112").Value = 16544.1
$ws.Range("N112").Value = -18760.1
# Row 115
$ws.Range("H115").Value = 844.5
$ws.Range("I115").Value = 190
$ws.Range("K115").Value = 570
$ws.Range("M115").Value = 997
# Row 122
$ws.Range("H122").Value = 1982.9756
$ws.Range("I122").Value = 1982.9756
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5948.9268
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3498.9268
$ws.Range("N122").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 4746.421
$ws.Range("I45").Value = 2551
$ws.Range("K45").Value = 2551
$ws.Range("M45").Value = -2174
# Row 63
$ws.Range("H63").Value = 2247.1667
$ws.Range("I63").Value = 2171.875
$ws.Range("J63").Value = 2397.75
$ws.Range("K63").Value = 2171.875
$ws.Range("L63").Value = 2397.75
$ws.Range("M63").Value = -1485.875
$ws.Range("N63").Value = -3769.75
# Row 66
$ws.Range("H66").Value = 2247.1667
$ws.Range("I66").Value = 2171.875
$ws.Range("J66").Value = 2397.75
$ws.Range("K66").Value = 10859.375
$ws.Range("L66").Value = 11988.75
$ws.Range("M66").Value = -7427.375
$ws.Range("N66").Value = -18852.75
# Row 110
$ws.Range("H110").Value = 55557140
$ws.Range("I110").Value = 1730.5
$ws.Range("J110").Value = 83334850
$ws.Range("K110").Value = 1730.5
$ws.Range("L110").Value = 83334850
$ws.Range("M110").Value = 314.5
$ws.Range("N110").Value = -83338940
# Row 132
$ws.Range("H132").Value = 4410.263
$ws.Range("I132").Value = 2144.025
$ws.Range("K132").Value = 6432.075000000001
$ws.Range("M132").Value = -3902.075000000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 6949706
$ws.Range("I20").Value = 8775898
$ws.Range("K20").Value = 8775898
$ws.Range("M20").Value = -8775651
# Row 94
$ws.Range("H94").Value = 2233.5217
$ws.Range("I94").Value = 1309.1875
$ws.Range("K94").Value = 1309.1875
$ws.Range("M94").Value = -858.1875
# Row 99
$ws.Range("H99").Value = 2068802
$ws.Range("I99").Value = 2425.0303
$ws.Range("K99").Value = 2425.0303
$ws.Range("M99").Value = -927.0302999999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7570.163
$ws.Range("I31").Value = 1798.7142
$ws.Range("K31").Value = 1798.7142
$ws.Range("M31").Value = -1503.7142
# Row 34
$ws.Range("H34").Value = 7570.163
$ws.Range("I34").Value = 1798.7142
$ws.Range("K34").Value = 1798.7142
$ws.Range("M34").Value = -1596.7142
# Row 62
$ws.Range("H62").Value = 5213493.5
$ws.Range("I62").Value = 12504320
$ws.Range("J62").Value = 5760.5713
$ws.Range("K62").Value = 12504320
$ws.Range("L62").Value = 5760.5713
$ws.Range("M62").Value = -12503696
$ws.Range("N62").Value = -7008.5713
# Row 65
$ws.Range("H65").Value = 5213493.5
$ws.Range("I65").Value = 12504320
$ws.Range("J65").Value = 5760.5713
$ws.Range("K65").Value = 62521600
$ws.Range("L65").Value = 28802.8565
$ws.Range("M65").Value = -62518480
$ws.Range("N65").Value = -35042.85649999999
# Row 100
$ws.Range("H100").Value = 45845.332
$ws.Range("J100").Value = 45845.332
$ws.Range("L100").Value = 45845.332
$ws.Range("N100").Value = -48009.332
# Row 107
$ws.Range("H107").Value = 2255.1177
$ws.Range("I107").Value = 790
$ws.Range("K107").Value = 790
$ws.Range("M107").Value = 1130
# Row 132
$ws.Range("H132").Value = 8256.291999999999
$ws.Range("I132").Value = 5709.4287
$ws.Range("K132").Value = 17128.2861
$ws.Range("M132").Value = -14598.2861
# Row 134
$ws.Range("H134").Value = 5129.0835
$ws.Range("I134").Value = 1764.2727
$ws.Range("J134").Value = 10416.643
$ws.Range("K134").Value = 5292.8181
$ws.Range("L134").Value = 31249.929
$ws.Range("M134").Value = -2757.8181
$ws.Range("N134").Value = -36319.929

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1751.05
$ws.Range("I5").Value = 601.4
$ws.Range("K5").Value = 1804.2
$ws.Range("M5").Value = -1692.2
# Row 37
$ws.Range("H37").Value = 66666
$ws.Range("J37").Value = 66666
$ws.Range("L37").Value = 199998
$ws.Range("N37").Value = -200222
# Row 135
$ws.Range("H135").Value = 1751.05
$ws.Range("I135").Value = 601.4
$ws.Range("K135").Value = 5412.599999999999
$ws.Range("M135").Value = -2877.599999999999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3588
$ws.Range("J80").Value = 4093.6
$ws.Range("L80").Value = 4093.6
$ws.Range("N80").Value = -6089.6
# Row 83
$ws.Range("H83").Value = 3588
$ws.Range("J83").Value = 4093.6
$ws.Range("L83").Value = 20468
$ws.Range("N83").Value = -30452
# Row 113
$ws.Range("H113").Value = 5090.567
$ws.Range("I113").Value = 3246.6365
$ws.Range("K113").Value = 3246.6365
$ws.Range("M113").Value = -1076.6365
# Row 132
$ws.Range("H132").Value = 4729.696
$ws.Range("I132").Value = 1924.25
$ws.Range("K132").Value = 5772.75
$ws.Range("M132").Value = -3242.75

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4194.6
$ws.Range("I7").Value = 3657.5293
$ws.Range("J7").Value = 5335.875
$ws.Range("K7").Value = 3657.5293
$ws.Range("L7").Value = 5335.875
$ws.Range("M7").Value = -3545.5293
$ws.Range("N7").Value = -5559.875
# Row 122
$ws.Range("H122").Value = 2553.6135
$ws.Range("I122").Value = 1804.6765
$ws.Range("K122").Value = 5414.029500000001
$ws.Range("M122").Value = -2964.029500000001
# Row 126
$ws.Range("H126").Value = 4194.6
$ws.Range("I126").Value = 3657.5293
$ws.Range("J126").Value = 5335.875
$ws.Range("K126").Value = 10972.5879
$ws.Range("L126").Value = 16007.625
$ws.Range("M126").Value = -8502.5879
$ws.Range("N126").Value = -20947.625

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 140963.8
$ws.Range("I122").Value = 268866.8
$ws.Range("K122").Value = 806600.3999999999
$ws.Range("M122").Value = -804150.3999999999
# Row 126
$ws.Range("H126").Value = 1051.8
$ws.Range("J126").Value = 1125
$ws.Range("L126").Value = 3375
$ws.Range("N126").Value = -8315
